{"js": "// Find every occurrence of the mis-ordered \"2022\" campaign-dates sentence\n// and reflow it so the year leads the sentence instead of trailing the\n// constellation name: \"Gwiazdozbi\u00f3r Oriona 2022:\" -> \"2022: ... Oriona:\".\nconst oldText = \": Daty kampanii u\u017cywaj\u0105ce Gwiazdozbi\u00f3r Oriona 2022: 16-25 stycznia, 14-23 lutego, 14-24 marca\";\nconst newText = \"2022: Daty kampanii u\u017cywaj\u0105ce Gwiazdozbi\u00f3r Oriona: 16-25 stycznia, 14-23 lutego, 14-24 marca\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Reflow the mis-ordered \"2022\" campaign-dates sentence so the year leads\n# the sentence instead of trailing the constellation name:\n# \"Gwiazdozbi\u00f3r Oriona 2022:\" -> \"2022: ... Oriona:\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \": Daty kampanii u\u017cywaj\u0105ce Gwiazdozbi\u00f3r Oriona 2022: 16-25 stycznia, 14-23 lutego, 14-24 marca\"\n$find.Replacement.Text = \"2022: Daty kampanii u\u017cywaj\u0105ce Gwiazdozbi\u00f3r Oriona: 16-25 stycznia, 14-23 lutego, 14-24 marca\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.Execute($find.Text, $find.MatchCase, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
